$d = $word.ActiveDocument
$d.Content.Find.Execute("Fogl Barna", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fogl Barna", 2)
